# Amazon Narrative.docx edit script
# 1) Expand "...utilised betweenness centrality to examine..." into
#    "...utilised betweenness centrality and key player analysis (see
#    Annex A) to examine..." with "key player analysis" underlined.
# 2) Change the References section (heading + 4 entries) from
#    double spacing (480) to 1.5 spacing (360).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: insert "and key player analysis (see Annex A) " before the
# existing "to examine the level of influence..." sentence.
# ---------------------------------------------------------------------

$curly_open  = [char]0x2018
$curly_close = [char]0x2019

$originalTail = " to examine the level of influence these Singaporean Parliamentarians have. Betweenness centrality measures the extent that a node sits " + $curly_open + "between" + $curly_close + " pairs of other nodes in the network. A node (i.e., MP) with high betweenness is prominent because that node is in a position to observe or control the flow of information in the network."

$target = $d.Content
$found = $target.Find.Execute($originalTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $newTail = " and key player analysis (see Annex A) to examine the level of influence these Singaporean Parliamentarians have. Betweenness centrality measures the extent that a node sits " + $curly_open + "between" + $curly_close + " pairs of other nodes in the network. A node (i.e., MP) with high betweenness is prominent because that node is in a position to observe or control the flow of information in the network. "

    # Re-writing .Text on the matched Range keeps the run's existing
    # rPr (Arial/Arial/Arial cs, sz 21 / szCs 21, no underline) intact
    # for the whole replacement instead of synthesising a bare <w:r>.
    $target.Text = $newTail

    # Now underline just the newly-added "key player analysis" phrase.
    $underlineRange = $d.Content
    $underlineRange.Find.Execute("key player analysis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $underlineRange.Font.Underline = 1
}

# ---------------------------------------------------------------------
# Part 2: References section line spacing 480 -> 360 (double -> 1.5).
# ---------------------------------------------------------------------

function Set-ParaLineSpacing15($searchText) {
    $c = $d.Content
    $c.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $p = $c.Paragraphs.First
    $p.Range.ParagraphFormat.LineSpacingRule = 1
}

Set-ParaLineSpacing15("References")
Set-ParaLineSpacing15("Campbell, A., Converse, P. E.")
Set-ParaLineSpacing15("Huddy, L., Sears, D. O.")
Set-ParaLineSpacing15("Martin, S. (2014). Parliamentary questions")
Set-ParaLineSpacing15("Staerklé, C. (2015)")

Write-Output "edit complete"
